# Fill in the first two empty rows of the time-tracking table with the
# new entries (commit: "Työaikoja ja muokkaa painike toimimaan").

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Find the first row whose "Pvm" (date) cell is still empty; the table's
# data rows are contiguous and end with a block of still-unused rows.
$targetRow = 0
for ($i = 2; $i -le $t.Rows.Count; $i++) {
    $cellText = $t.Cell($i, 1).Range.Text
    # An empty cell's Range.Text is just the cell-end marker (length 2).
    if ($cellText.Length -le 2) {
        $targetRow = $i
        break
    }
}

$t.Cell($targetRow, 1).Range.Text = "03.03.2023"
$t.Cell($targetRow, 2).Range.Text = "0,5h"
$t.Cell($targetRow, 3).Range.Text = "Sprintti palaveri"

$t.Cell($targetRow + 1, 1).Range.Text = "06.03.2023"
$t.Cell($targetRow + 1, 2).Range.Text = "0,5h"
$t.Cell($targetRow + 1, 3).Range.Text = "Muokkaa painike toimimaan"
